$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.500.57"
$ws.Range("E2").Value = "  +0.42%  "

# Row 3
$ws.Range("D3").Value = "3.444.82"
$ws.Range("E3").Value = "  +1.21%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.63"
$ws.Range("E5").Value = "  +0.70%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.17"
$ws.Range("E6").Value = "  +4.54%  "

# Row 7
$ws.Range("D7").Value = "3.445.13"
$ws.Range("E7").Value = "  +1.24%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.478"
$ws.Range("E9").Value = "  +2.18%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.69"
$ws.Range("E10").Value = "  +0.15%  "

# Row 11
$ws.Range("E11").Value = "  +3.22%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.389"
$ws.Range("E12").Value = "  +2.19%  "

# Row 13
$ws.Range("D13").Value = "4.031.24"
$ws.Range("E13").Value = "  +1.23%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.49"
$ws.Range("E14").Value = "  +6.77%  "

# Row 15
$ws.Range("E15").Value = "  -0.33%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000174"
$ws.Range("E16").Value = "  +0.91%  "

# Row 17
$ws.Range("D17").Value = "3.445.84"
$ws.Range("E17").Value = "  +1.33%  "

# Row 18
$ws.Range("D18").Value = "61.624.03"
$ws.Range("E18").Value = "  +0.60%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.34"
$ws.Range("E19").Value = "  +7.09%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.37"
$ws.Range("E20").Value = "  +3.59%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.45"
$ws.Range("E21").Value = "  +1.03%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "403.06"
$ws.Range("E22").Value = "  +7.03%  "

# Row 23
$ws.Range("E23").Value = "  +3.02%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.05"
$ws.Range("E24").Value = "  +4.03%  "

# Row 25
$ws.Range("E25").Value = "  +0.51%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.15%  "

# Row 27
$ws.Range("E27").Value = "  -0.82%  "

# Row 28
$ws.Range("D28").Value = "3.584.20"
$ws.Range("E28").Value = "  +1.66%  "

# Row 29
$ws.Range("E29").Value = "  +4.52%  "

# Row 30
$ws.Range("E30").Value = "  +2.41%  "

# Row 31
$ws.Range("E31").Value = "  +0.02%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.29"
$ws.Range("E32").Value = "  +2.45%  "

# Row 33
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.19"
$ws.Range("E33").Value = "  +2.08%  "

# Row 34
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.46"
$ws.Range("E34").Value = "  -10.54%  "

# Row 35
$ws.Range("E35").Value = "  -0.05%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.94"
$ws.Range("E36").Value = "  +1.99%  "

# Row 37
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.04"
$ws.Range("E37").Value = "  +2.46%  "

# Row 38
$ws.Range("B38").Value = "RenzoRestakedETH"
$ws.Range("C38").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D38").Value = "3.472.40"
$ws.Range("E38").Value = "  +1.44%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.15"
$ws.Range("E39").Value = "  +0.36%  "

# Row 40
$ws.Range("E40").Value = "  -0.26%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "167.18"
$ws.Range("E41").Value = "  +0.84%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0792"
$ws.Range("E42").Value = "  +2.77%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "27.24"
$ws.Range("E43").Value = "  +3.83%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.803"
$ws.Range("E44").Value = "  +3.30%  "

# Row 45
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.53"
$ws.Range("E45").Value = "  +3.05%  "

# Row 46
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.05%  "

# Row 47
$ws.Range("E47").Value = "  -1.95%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "42.40"
$ws.Range("E48").Value = "  +1.04%  "

# Row 49
$ws.Range("D49").Value = "2.607.98"
$ws.Range("E49").Value = "  +2.89%  "

# Row 50
$ws.Range("E50").Value = "  -2.14%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.96"
$ws.Range("E51").Value = "  +2.44%  "
